$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "42.071.57"
$r.Style = "Normal"
$r = $ws.Range("E2")
$r.NumberFormat = "@"
$r.Value = "  -3.68%  "
$r.Style = "Normal"

$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "2.204.59"
$r.Style = "Normal"
$r = $ws.Range("E3")
$r.NumberFormat = "@"
$r.Value = "  -3.22%  "
$r.Style = "Normal"

$r = $ws.Range("E4")
$r.NumberFormat = "@"
$r.Value = "  +0.38%  "
$r.Style = "Normal"

$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "106.80"
$r.Style = "Normal"
$r = $ws.Range("E5")
$r.NumberFormat = "@"
$r.Value = "  -14.17%  "
$r.Style = "Normal"

$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "293.04"
$r.Style = "Normal"
$r = $ws.Range("E6")
$r.NumberFormat = "@"
$r.Value = "  +9.89%  "
$r.Style = "Normal"

$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "0.619"
$r.Style = "Normal"
$r = $ws.Range("E7")
$r.NumberFormat = "@"
$r.Value = "  -3.24%  "
$r.Style = "Normal"

$r = $ws.Range("E8")
$r.NumberFormat = "@"
$r.Value = "  -0.10%  "
$r.Style = "Normal"

$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.591"
$r.Style = "Normal"
$r = $ws.Range("E9")
$r.NumberFormat = "@"
$r.Value = "  -5.63%  "
$r.Style = "Normal"

$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "43.63"
$r.Style = "Normal"
$r = $ws.Range("E10")
$r.NumberFormat = "@"
$r.Value = "  -9.64%  "
$r.Style = "Normal"

$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.0904"
$r.Style = "Normal"
$r = $ws.Range("E11")
$r.NumberFormat = "@"
$r.Value = "  -4.45%  "
$r.Style = "Normal"

$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "53.89"
$r.Style = "Normal"
$r = $ws.Range("E12")
$r.NumberFormat = "@"
$r.Value = "  -0.94%  "
$r.Style = "Normal"

$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "8.70"
$r.Style = "Normal"
$r = $ws.Range("E13")
$r.NumberFormat = "@"
$r.Value = "  -6.19%  "
$r.Style = "Normal"

$r = $ws.Range("E14")
$r.NumberFormat = "@"
$r.Value = "  -3.39%  "
$r.Style = "Normal"

$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "0.939"
$r.Style = "Normal"
$r = $ws.Range("E15")
$r.NumberFormat = "@"
$r.Value = "  +3.91%  "
$r.Style = "Normal"

$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "14.80"
$r.Style = "Normal"
$r = $ws.Range("E16")
$r.NumberFormat = "@"
$r.Value = "  -4.27%  "
$r.Style = "Normal"

$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "2.540.82"
$r.Style = "Normal"
$r = $ws.Range("E17")
$r.NumberFormat = "@"
$r.Value = "  -3.22%  "
$r.Style = "Normal"

$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "2.234.80"
$r.Style = "Normal"
$r = $ws.Range("E18")
$r.NumberFormat = "@"
$r.Value = "  -1.85%  "
$r.Style = "Normal"

$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "42.195.43"
$r.Style = "Normal"
$r = $ws.Range("E19")
$r.NumberFormat = "@"
$r.Value = "  -3.47%  "
$r.Style = "Normal"

$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "7.19"
$r.Style = "Normal"
$r = $ws.Range("E20")
$r.NumberFormat = "@"
$r.Value = "  +2.86%  "
$r.Style = "Normal"

$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "0.0000104"
$r.Style = "Normal"
$r = $ws.Range("E21")
$r.NumberFormat = "@"
$r.Value = "  -5.71%  "
$r.Style = "Normal"

$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "72.23"
$r.Style = "Normal"
$r = $ws.Range("E22")
$r.NumberFormat = "@"
$r.Value = "  -0.28%  "
$r.Style = "Normal"

$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "3.43"
$r.Style = "Normal"
$r = $ws.Range("E23")
$r.NumberFormat = "@"
$r.Value = "  +18.50%  "
$r.Style = "Normal"

$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "2.26"
$r.Style = "Normal"
$r = $ws.Range("E24")
$r.NumberFormat = "@"
$r.Value = "  -7.35%  "
$r.Style = "Normal"

$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "226.13"
$r.Style = "Normal"
$r = $ws.Range("E25")
$r.NumberFormat = "@"
$r.Value = "  -4.14%  "
$r.Style = "Normal"

$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "8.89"
$r.Style = "Normal"
$r = $ws.Range("E26")
$r.NumberFormat = "@"
$r.Value = "  -6.49%  "
$r.Style = "Normal"

$r = $ws.Range("E27")
$r.NumberFormat = "@"
$r.Value = "  -1.35%  "
$r.Style = "Normal"

$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "11.50"
$r.Style = "Normal"
$r = $ws.Range("E28")
$r.NumberFormat = "@"
$r.Value = "  -3.68%  "
$r.Style = "Normal"

$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "3.93"
$r.Style = "Normal"
$r = $ws.Range("E29")
$r.NumberFormat = "@"
$r.Value = "  +0.31%  "
$r.Style = "Normal"

$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "2.23"
$r.Style = "Normal"
$r = $ws.Range("E30")
$r.NumberFormat = "@"
$r.Value = "  -1.56%  "
$r.Style = "Normal"

$r = $ws.Range("B31")
$r.NumberFormat = "@"
$r.Value = "InjectiveProtocol"
$r.Style = "Normal"
$r = $ws.Range("C31")
$r.NumberFormat = "@"
$r.Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$r.Style = "Normal"
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "37.79"
$r.Style = "Normal"
$r = $ws.Range("E31")
$r.NumberFormat = "@"
$r.Value = "  -11.51%  "
$r.Style = "Normal"

$r = $ws.Range("B32")
$r.NumberFormat = "@"
$r.Value = "WEMIXToken"
$r.Style = "Normal"
$r = $ws.Range("C32")
$r.NumberFormat = "@"
$r.Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$r.Style = "Normal"
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "3.20"
$r.Style = "Normal"
$r = $ws.Range("E32")
$r.NumberFormat = "@"
$r.Value = "  -4.90%  "
$r.Style = "Normal"

$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "172.48"
$r.Style = "Normal"
$r = $ws.Range("E33")
$r.NumberFormat = "@"
$r.Value = "  -0.21%  "
$r.Style = "Normal"

$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "20.72"
$r.Style = "Normal"
$r = $ws.Range("E34")
$r.NumberFormat = "@"
$r.Value = "  -4.60%  "
$r.Style = "Normal"

$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "0.0873"
$r.Style = "Normal"
$r = $ws.Range("E35")
$r.NumberFormat = "@"
$r.Value = "  -5.35%  "
$r.Style = "Normal"

$r = $ws.Range("B36")
$r.NumberFormat = "@"
$r.Value = "RenderToken"
$r.Style = "Normal"
$r = $ws.Range("C36")
$r.NumberFormat = "@"
$r.Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$r.Style = "Normal"
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "5.02"
$r.Style = "Normal"
$r = $ws.Range("E36")
$r.NumberFormat = "@"
$r.Value = "  +7.81%  "
$r.Style = "Normal"

$r = $ws.Range("B37")
$r.NumberFormat = "@"
$r.Value = "Filecoin"
$r.Style = "Normal"
$r = $ws.Range("C37")
$r.NumberFormat = "@"
$r.Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$r.Style = "Normal"
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "5.47"
$r.Style = "Normal"
$r = $ws.Range("E37")
$r.NumberFormat = "@"
$r.Value = "  -5.32%  "
$r.Style = "Normal"

$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "4.26"
$r.Style = "Normal"
$r = $ws.Range("E38")
$r.NumberFormat = "@"
$r.Value = "  -0.50%  "
$r.Style = "Normal"

$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "0.125"
$r.Style = "Normal"
$r = $ws.Range("E39")
$r.NumberFormat = "@"
$r.Value = "  -4.04%  "
$r.Style = "Normal"

$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "0.0359"
$r.Style = "Normal"
$r = $ws.Range("E40")
$r.NumberFormat = "@"
$r.Value = "  -5.26%  "
$r.Style = "Normal"

$r = $ws.Range("E41")
$r.NumberFormat = "@"
$r.Value = "  -4.98%  "
$r.Style = "Normal"

$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "2.44"
$r.Style = "Normal"
$r = $ws.Range("E42")
$r.NumberFormat = "@"
$r.Value = "  -4.46%  "
$r.Style = "Normal"

$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "69.68"
$r.Style = "Normal"
$r = $ws.Range("E43")
$r.NumberFormat = "@"
$r.Value = "  -6.36%  "
$r.Style = "Normal"

$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "0.227"
$r.Style = "Normal"
$r = $ws.Range("E44")
$r.NumberFormat = "@"
$r.Value = "  -5.02%  "
$r.Style = "Normal"

$r = $ws.Range("E45")
$r.NumberFormat = "@"
$r.Value = "  +0.15%  "
$r.Style = "Normal"

$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "12.55"
$r.Style = "Normal"
$r = $ws.Range("E46")
$r.NumberFormat = "@"
$r.Value = "  -10.12%  "
$r.Style = "Normal"

$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "1.28"
$r.Style = "Normal"
$r = $ws.Range("E47")
$r.NumberFormat = "@"
$r.Value = "  -6.48%  "
$r.Style = "Normal"

$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "5.37"
$r.Style = "Normal"
$r = $ws.Range("E48")
$r.NumberFormat = "@"
$r.Value = "  -4.64%  "
$r.Style = "Normal"

$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "1.30"
$r.Style = "Normal"
$r = $ws.Range("E49")
$r.NumberFormat = "@"
$r.Value = "  +2.70%  "
$r.Style = "Normal"

$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "101.83"
$r.Style = "Normal"
$r = $ws.Range("E50")
$r.NumberFormat = "@"
$r.Value = "  -0.23%  "
$r.Style = "Normal"

$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "8.36"
$r.Style = "Normal"
$r = $ws.Range("E51")
$r.NumberFormat = "@"
$r.Value = "  -2.90%  "
$r.Style = "Normal"

